$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.318.24"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").Value = "3.429.92"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.49"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.86"
$ws.Range("E6").Value = "  +3.65%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "8.07"
$ws.Range("E9").Value = "  +4.88%  "
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("D12").Value = "4.018.45"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.67"
$ws.Range("E14").Value = "  -3.41%  "
$ws.Range("D15").Value = "3.438.76"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "62.333.81"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.53"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.58"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.99"
$ws.Range("E20").Value = "  -3.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.73"
$ws.Range("E21").Value = "  -1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.572"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "75.31"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "3.562.31"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.181"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.71"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.99"
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.27"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.62"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.95"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "31.10"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").Value = "3.463.98"
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0789"
$ws.Range("E41").Value = "  +3.00%  "
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.779"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.41"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").Value = "2.548.19"
$ws.Range("E47").Value = "  -2.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.89"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("E49").Value = "  -4.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.66"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("E51").Value = "  -0.09%  "
